# Commit "Add files via upload" — the uploaded replacement workbook
# redacts/anonymizes a few personal-data strings on the "tweets" sheet
# (row 5 col A, row 6 col B, row 7 col A) and replaces the hashtag/mention
# list on the "hashtag" sheet with placeholder text. The "say" sheet's
# single numeric cell is untouched. The previously-active sheet ("say")
# is no longer selected; "tweets" becomes the active sheet/tab instead.

$wb = $excel.ActiveWorkbook

$tweets  = $wb.Worksheets.Item("tweets")
$hashtag = $wb.Worksheets.Item("hashtag")
$say     = $wb.Worksheets.Item("say")

# --- Redact / replace specific cell contents on "tweets" ---
# Row 5, col A: "KONUTZEDELER çözüm bekliyor." -> "Çözüm bekliyoruz"
$tweets.Range("A5").Value = "Çözüm bekliyoruz"

# Row 6, col B: "Biz bu evleri almak için hayatımızı ortaya koyduk" -> redacted
$tweets.Range("B6").Value = "Biz bu …. almak için hayatımızı ortaya koyduk"

# Row 7, col A: "biz kendimizi konutzede değil ‘terörzede’ olarak görüyoruz." -> redacted
$tweets.Range("A7").Value = "biz kendimiz i…. değil ‘terörzede’ olarak görüyoruz."

# --- Replace the hashtag/mentions cell on "hashtag" with placeholders ---
$hashtag.Range("A1").Value = "#hashtag1 #hashtag2 #hashtag3 @kişi1 @kişi2"

# --- "say" sheet numeric cell is unchanged ---
# ($say.Range("A1").Value stays 7)

# --- Update selection / active sheet state ---
$hashtag.Range("A2").Select() | Out-Null
$say.Range("A1").Select() | Out-Null
$tweets.Activate() | Out-Null
$tweets.Range("A9").Select() | Out-Null
